$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row before the "Description" row (currently row 11),
# which pushes Description..Context down by one row.
$ws.Rows.Item(11).Insert()

# Populate the new "Jurisdiction" property row (value left blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the "Date" property value.
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
